$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new values
$ws.Range("A2").Value = 45113.50694444445
$ws.Range("B2").Value = 10.726
$ws.Range("C2").Value = 7.333
$ws.Range("D2").Value = 3.404
$ws.Range("E2").Value = 23.56
$ws.Range("F2").Value = 17.15
$ws.Range("G2").Value = 8.176
$ws.Range("H2").Value = 24.228
$ws.Range("I2").Value = 13.347
$ws.Range("J2").Value = 5.245
$ws.Range("K2").Value = 7.323
$ws.Range("L2").Value = 9.308999999999999
$ws.Range("M2").Value = 10.191
$ws.Range("N2").Value = 2.44
$ws.Range("O2").Value = 8.647
$ws.Range("P2").Value = 11.655
$ws.Range("Q2").Value = 7.955
$ws.Range("R2").Value = 2.648
$ws.Range("S2").Value = 1.093
$ws.Range("T2").Value = 124.223
$ws.Range("U2").Value = 23.834
$ws.Range("V2").Value = 7.982
$ws.Range("W2").Value = 14.964
$ws.Range("X2").Value = 8.048999999999999
$ws.Range("Y2").Value = 2.19
$ws.Range("Z2").Value = 13.597
$ws.Range("AA2").Value = 7.05
$ws.Range("AB2").Value = 6.629
$ws.Range("AC2").Value = 7.562
$ws.Range("AD2").Value = 9.942
$ws.Range("AE2").Value = 2.682
$ws.Range("AF2").Value = 21.628
$ws.Range("AG2").Value = 4.075
$ws.Range("AH2").Value = 9.978

$ws.Range("A3").Value = 45113.51388888889
$ws.Range("B3").Value = 17.078
$ws.Range("C3").Value = 12.559
$ws.Range("D3").Value = 1.781
$ws.Range("E3").Value = 37.536
$ws.Range("F3").Value = 29.804
$ws.Range("G3").Value = 13.317
$ws.Range("H3").Value = 49.837
$ws.Range("I3").Value = 20.911
$ws.Range("J3").Value = 9.137
$ws.Range("K3").Value = 13.127
$ws.Range("L3").Value = 15.023
$ws.Range("M3").Value = 16.115
$ws.Range("N3").Value = 4.126
$ws.Range("O3").Value = 13.535
$ws.Range("P3").Value = 19.015
$ws.Range("Q3").Value = 11.748
$ws.Range("R3").Value = 1.343
$ws.Range("S3").Value = 0.795
$ws.Range("T3").Value = 198.684
$ws.Range("U3").Value = 37.811
$ws.Range("V3").Value = 12.493
$ws.Range("W3").Value = 25.002
$ws.Range("X3").Value = 13.163
$ws.Range("Y3").Value = 2.184
$ws.Range("Z3").Value = 25.307
$ws.Range("AA3").Value = 11.035
$ws.Range("AB3").Value = 9.955
$ws.Range("AC3").Value = 11.647
$ws.Range("AD3").Value = 15.816
$ws.Range("AE3").Value = 1.136
$ws.Range("AF3").Value = 45.487
$ws.Range("AG3").Value = 6.848
$ws.Range("AH3").Value = 15.619

$ws.Range("A4").Value = 45113.52083333334
$ws.Range("B4").Value = 9.445
$ws.Range("C4").Value = 6.946
$ws.Range("D4").Value = 1.071
$ws.Range("E4").Value = 20.866
$ws.Range("F4").Value = 16.415
$ws.Range("G4").Value = 7.357
$ws.Range("H4").Value = 31.874
$ws.Range("I4").Value = 11.605
$ws.Range("J4").Value = 5.099
$ws.Range("K4").Value = 7.16
$ws.Range("L4").Value = 8.356
$ws.Range("M4").Value = 9.006
$ws.Range("N4").Value = 2.252
$ws.Range("O4").Value = 7.519
$ws.Range("P4").Value = 10.545
$ws.Range("Q4").Value = 6.615
$ws.Range("R4").Value = 0.878
$ws.Range("S4").Value = 0.445
$ws.Range("T4").Value = 107.138
$ws.Range("U4").Value = 21.135
$ws.Range("V4").Value = 6.941
$ws.Range("W4").Value = 13.905
$ws.Range("X4").Value = 7.286
$ws.Range("Y4").Value = 1.239
$ws.Range("Z4").Value = 15.508
$ws.Range("AA4").Value = 6.131
$ws.Range("AB4").Value = 5.574
$ws.Range("AC4").Value = 6.518
$ws.Range("AD4").Value = 8.792
$ws.Range("AE4").Value = 0.722
$ws.Range("AF4").Value = 29.214
$ws.Range("AG4").Value = 3.763
$ws.Range("AH4").Value = 8.678000000000001

$ws.Range("A5").Value = 45113.52777777778
$ws.Range("B5").Value = 17.64
$ws.Range("C5").Value = 13.16
$ws.Range("D5").Value = 1.14
$ws.Range("E5").Value = 38.65
$ws.Range("F5").Value = 31.32
$ws.Range("G5").Value = 13.83
$ws.Range("H5").Value = 52.64
$ws.Range("I5").Value = 21.5
$ws.Range("J5").Value = 9.58
$ws.Range("K5").Value = 13.94
$ws.Range("L5").Value = 15.5
$ws.Range("M5").Value = 16.53
$ws.Range("N5").Value = 4.34
$ws.Range("O5").Value = 13.91
$ws.Range("P5").Value = 19.71
$ws.Range("Q5").Value = 11.82
$ws.Range("R5").Value = 0.75
$ws.Range("S5").Value = 0.63
$ws.Range("T5").Value = 204.41
$ws.Range("U5").Value = 38.85
$ws.Range("V5").Value = 12.84
$ws.Range("W5").Value = 26.01
$ws.Range("X5").Value = 13.67
$ws.Range("Y5").Value = 2.01
$ws.Range("Z5").Value = 26.03
$ws.Range("AA5").Value = 11.34
$ws.Range("AB5").Value = 10.1
$ws.Range("AC5").Value = 11.86
$ws.Range("AD5").Value = 16.3
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 47.72
$ws.Range("AG5").Value = 7.17
$ws.Range("AH5").Value = 16.05

# Remove the now-unused row 6 (data trimmed from 5 to 4 rows)
$ws.Rows.Item(6).Delete()

# Adjust column widths (stored OOXML width = ColumnWidth + 5/6)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 6.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(14).ColumnWidth = 6.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(18).ColumnWidth = 6.166666666666667
$ws.Columns.Item(19).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(25).ColumnWidth = 6.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 6.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(31).ColumnWidth = 6.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(33).ColumnWidth = 6.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
